$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing columns (maxTds/minTds leftovers + fertilizer/nutrient
# pump telemetry) that no longer belong in the sample export.
$ws.Range("G1:L11").EntireColumn.Delete()

# Headers (row 1) for the remaining columns
$headers = @("datetime", "airTemperature", "humidity", "waterTemperature", "ph", "tds")
for ($col = 1; $col -le 6; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Refreshed data rows 2-11: datetime, airTemperature, humidity, waterTemperature, ph, tds
$data = @(
    @(45598.04027777778, 29.39999999999999, 84.00857142857141, 26.69000000000002, 6.360857142857142, 882.8571428571429),
    @(45598.17222222222, 33.08857142857143, 68.68333333333338, 26.88, 6.413888888888888, 876.4444444444445),
    @(45609.45486111111, 29.5, 84.67428571428574, 28.5, 6.286857142857143, 899.5428571428571),
    @(45613.02777777778, 29.79969135802468, 84.94857142857146, 26.99828571428572, 6.322571428571426, 812.2857142857143),
    @(45621.03958333333, 29.41388888888888, 85.75833333333334, 27.05999999999997, 6.335277777777778, 750.9166666666666),
    @(45622.27013888889, 32.05428571428572, 77.72571428571429, 27.25, 6.353714285714286, 735.8),
    @(45624.02777777778, 29.8742857142857, 88.32000000000004, 26.79114285714284, 6.321714285714286, 743.3428571428572),
    @(45625.65902777778, 27, 97.89444444444452, 27.38, 6.320277777777776, 775.5),
    @(45627.125, 26.95428571428571, 97.92857142857146, 26.44000000000002, 6.239428571428573, 987.5714285714286),
    @(45627.64444444444, 25.82857142857141, 100, 26.54742857142855, 6.259428571428571, 975.3714285714286)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
